$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '257.81'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.07%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.47'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.83%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.568'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-12.64%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05891'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.51%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.633'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.86%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8592'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-1.25%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9252'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-12.33%'

$ws.Range("B9").Value = 'One'

$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0006048'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.83%'

$ws.Range("B10").Value = 'WazirX'

$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1411'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.33%'

$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03659'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.08%'

$ws.Range("B12").Value = 'MandalaExchangeToken'

$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07088'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.61%'

$ws.Range("B13").Value = 'BitrueCoin'

$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03173'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-2.10%'

$ws.Range("B14").Value = 'BitMartToken'

$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09178'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.43%'

$ws.Range("B15").Value = 'BitForexToken'

$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001538'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.27%'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006090'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.63%'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.514'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.85%'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.198'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.17%'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3106'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.38%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1278'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.13%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.863'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '8.80%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04213'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.50%'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.41%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004302'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-5.31%'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.15%'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-22.20%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03838'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.38%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006248'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '56.55%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1101'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.51%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002199'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-10.71%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01146'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '15.22%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005457'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.38%'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.18%'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-45.10%'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1340'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '6,160.09%'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.18%'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.18%'
